# Add a new machine record (row 33) to the master-machine_master sheet,
# and update the sheet's view/selection state to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (mirrors the pattern of the existing rows: id, name,
# mac_address, serial_num, ip_address, mspec_id, lang_code, is_active,
# cr_by, cr_dtimes, eff_dtimes)
$ws.Range("A33").Value = 10032
$ws.Range("B33").Value = "Machine 32"
$ws.Range("C33").Value = "F4-30-B9-D4-CD-6F"
$ws.Range("D33").Value = "FB5962911665"
$ws.Range("E33").Value = "192.168.0.358"
$ws.Range("F33").Value = 1001
$ws.Range("G33").Value = "eng"
$ws.Range("H33").Value = $true
$ws.Range("I33").Value = "superadmin"
$ws.Range("J33").Value = "now()"
$ws.Range("K33").Value = "now()"

# Scroll the view down toward the newly added row and select it.
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C28").Select() | Out-Null
